$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet
$ws.Name = "IngestMetadataFile"

# Update the header row (row 1) values
$headers = @(
    "file_name",
    "file_created_by",
    "file_creation_date",
    "ingest_code_url",
    "source_infores_id",
    "source_data_version",
    "source_access_date",
    "source_access_urls",
    "source_file_names",
    "target_name",
    "target_creation_date",
    "target_format",
    "target_model",
    "target_model_url",
    "target_data_model_version",
    "node_normalizer",
    "node_normalizer_version",
    "node_normalizer_url",
    "total_edge_count",
    "total_node_count",
    "orphan_node_count",
    "node_categories",
    "edge_predicates"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
